$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

# Copy formatting from the last existing row (610) so new rows match style s="1"
$ws.Range("A610:C610").Copy() | Out-Null
$ws.Range("A611:C625").PasteSpecial(-4122) | Out-Null

$ws.Range("A611").Value = "cs"
$ws.Range("B611").Value = "common.taste.apple"
$ws.Range("C611").Value = "Jablko"
$ws.Range("A612").Value = "cs"
$ws.Range("B612").Value = "common.taste.pear"
$ws.Range("C612").Value = "Hruška"
$ws.Range("A613").Value = "cs"
$ws.Range("B613").Value = "common.taste.apricot"
$ws.Range("C613").Value = "Meruňka"
$ws.Range("A614").Value = "cs"
$ws.Range("B614").Value = "common.taste.cinnamon"
$ws.Range("C614").Value = "Skořice"
$ws.Range("A615").Value = "cs"
$ws.Range("B615").Value = "common.taste.mint"
$ws.Range("C615").Value = "Máta"
$ws.Range("A616").Value = "cs"
$ws.Range("B616").Value = "common.taste.grape"
$ws.Range("C616").Value = "Hrozen"
$ws.Range("A617").Value = "cs"
$ws.Range("B617").Value = "common.taste.hazelnut"
$ws.Range("C617").Value = "Oříšky"
$ws.Range("A618").Value = "cs"
$ws.Range("B618").Value = "common.taste.cherry"
$ws.Range("C618").Value = "Třešeň"
$ws.Range("A619").Value = "cs"
$ws.Range("B619").Value = "common.taste.bourbon"
$ws.Range("C619").Value = "Bourbon"
$ws.Range("A620").Value = "cs"
$ws.Range("B620").Value = "common.taste.orange"
$ws.Range("C620").Value = "Pomeranč"
$ws.Range("A621").Value = "cs"
$ws.Range("B621").Value = "common.taste.cookie"
$ws.Range("C621").Value = "Sušenka"
$ws.Range("A622").Value = "cs"
$ws.Range("B622").Value = "common.taste.lemon"
$ws.Range("C622").Value = "Citron"
$ws.Range("A623").Value = "cs"
$ws.Range("B623").Value = "common.taste.peanut"
$ws.Range("C623").Value = "Arašídy"
$ws.Range("A624").Value = "cs"
$ws.Range("B624").Value = "common.taste.pistachio"
$ws.Range("C624").Value = "Pistácie"
$ws.Range("A625").Value = "cs"
$ws.Range("B625").Value = "common.taste.plum"
$ws.Range("C625").Value = "Švestka"

# Update selection / active cell on the Translations - Common sheet and activate it
$ws.Activate() | Out-Null
$ws.Range("B620").Select() | Out-Null
